# Fruta / hortaliza, semanal
# Insert two new weekly price records for Piña (Vega Modelo de Temuco) at
# the top of the data block (rows 457-458), pushing every existing record
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 457, shifting rows 457:489 down to 459:491.
$ws.Range("A457:A458").EntireRow.Insert()

# Row 457: new "Primera" quality record.
$ws.Cells.Item(457, 1).Value = 10
$ws.Cells.Item(457, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(457, 3).Value = "La Araucanía"
$ws.Cells.Item(457, 4).Value = 44746
$ws.Cells.Item(457, 5).Value = 9
$ws.Cells.Item(457, 6).Value = "Fruta"
$ws.Cells.Item(457, 7).Value = 100108
$ws.Cells.Item(457, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(457, 9).Value = 100108005
$ws.Cells.Item(457, 10).Value = "Piña"
$ws.Cells.Item(457, 11).Value = "Caramelo"
$ws.Cells.Item(457, 12).Value = "Primera"
$ws.Cells.Item(457, 13).Value = 90
$ws.Cells.Item(457, 14).Value = 24000
$ws.Cells.Item(457, 15).Value = 24000
$ws.Cells.Item(457, 16).Value = 24000
$ws.Cells.Item(457, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(457, 18).Value = "Ecuador"
$ws.Cells.Item(457, 19).Value = 2000
$ws.Cells.Item(457, 20).Value = 12

# Row 458: new "Segunda" quality record.
$ws.Cells.Item(458, 1).Value = 10
$ws.Cells.Item(458, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(458, 3).Value = "La Araucanía"
$ws.Cells.Item(458, 4).Value = 44746
$ws.Cells.Item(458, 5).Value = 9
$ws.Cells.Item(458, 6).Value = "Fruta"
$ws.Cells.Item(458, 7).Value = 100108
$ws.Cells.Item(458, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(458, 9).Value = 100108005
$ws.Cells.Item(458, 10).Value = "Piña"
$ws.Cells.Item(458, 11).Value = "Caramelo"
$ws.Cells.Item(458, 12).Value = "Segunda"
$ws.Cells.Item(458, 13).Value = 220
$ws.Cells.Item(458, 14).Value = 21000
$ws.Cells.Item(458, 15).Value = 23000
$ws.Cells.Item(458, 16).Value = 21909
$ws.Cells.Item(458, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(458, 18).Value = "Ecuador"
$ws.Cells.Item(458, 19).Value = 1565
$ws.Cells.Item(458, 20).Value = 14
